$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J16").Value = 0.2150495036779461
$ws.Range("I17").Value = 0.24
$ws.Range("H18").Value = 0.3087982760018804
$ws.Range("G19").Value = 0.32
$ws.Range("F20").Value = 0.4476495795507702
$ws.Range("E21").Value = 0.1088966743764388
$ws.Range("D22").Value = 0.1461563307127136
$ws.Range("C23").Value = 0.09547648014918764
$ws.Range("B24").Value = 0.0959495356205764
